# Update tracker data 2025-09-18 00:02:05
# Appends one new day's worth of progress-tracker rows (date serial 45918,
# i.e. 2025-09-18) for each of the 5 goals, mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45918
$progress = 0.896323717517805
$percentage = 0
$change = -0.01

$goals = @(
    @("G2", "Workout"),
    @("G3", "Eat Healthy"),
    @("G4", "Read Book"),
    @("G5", "Investment Plan"),
    @("G6", "Spend 10 Hours without phone")
)

$startRow = 57
for ($i = 0; $i -lt $goals.Count; $i++) {
    $row = $startRow + $i
    $goalId = $goals[$i][0]
    $goalName = $goals[$i][1]

    $ws.Cells.Item($row, 1).Value = $goalId
    $ws.Cells.Item($row, 2).Value = $goalName

    $ws.Cells.Item($row, 3).Value = $newDate
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = $progress
    $ws.Cells.Item($row, 5).Value = $percentage
    $ws.Cells.Item($row, 6).Value = $change
}
